$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its text formatting so numeric-looking
# values like "209.19" are not auto-converted into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.122.86"
$ws.Range("E2").Value = "  -2.26%  "

$ws.Range("D3").Value = "1.576.11"
$ws.Range("E3").Value = "  -1.57%  "

$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").Value = "209.19"
$ws.Range("E5").Value = "  -1.10%  "

$ws.Range("D6").Value = "0.496"
$ws.Range("E6").Value = "  -3.44%  "

$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("E8").Value = "  -0.70%  "

$ws.Range("D9").Value = "0.0608"
$ws.Range("E9").Value = "  -1.63%  "

$ws.Range("D10").Value = "19.52"
$ws.Range("E10").Value = "  -0.86%  "

$ws.Range("D11").Value = "0.0844"
$ws.Range("E11").Value = "  -0.24%  "

$ws.Range("D12").Value = "1.798.94"
$ws.Range("E12").Value = "  -1.45%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.580.66"
$ws.Range("E13").Value = "  -1.14%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  -0.19%  "

$ws.Range("D15").Value = "0.512"
$ws.Range("E15").Value = "  -2.15%  "

$ws.Range("D16").Value = "64.36"
$ws.Range("E16").Value = "  -1.03%  "

$ws.Range("D17").Value = "26.137.39"
$ws.Range("E17").Value = "  -2.07%  "

$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  -1.53%  "

$ws.Range("E19").Value = "  +1.43%  "

$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("D21").Value = "206.96"
$ws.Range("E21").Value = "  -1.35%  "

$ws.Range("D22").Value = "4.24"
$ws.Range("E22").Value = "  -0.97%  "

$ws.Range("E23").Value = "  -1.63%  "

$ws.Range("D24").Value = "8.86"
$ws.Range("E24").Value = "  -1.29%  "

$ws.Range("D25").Value = "144.21"
$ws.Range("E25").Value = "  +0.52%  "

$ws.Range("E26").Value = "  -0.27%  "

$ws.Range("D27").Value = "6.97"
$ws.Range("E27").Value = "  -1.53%  "

$ws.Range("E28").Value = "  -2.01%  "

$ws.Range("D29").Value = "15.21"
$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("E31").Value = "  -1.14%  "

$ws.Range("E32").Value = "  -1.96%  "

$ws.Range("D33").Value = "2.97"
$ws.Range("E33").Value = "  +0.18%  "

$ws.Range("D34").Value = "1.279.35"
$ws.Range("E34").Value = "  -0.77%  "

$ws.Range("E35").Value = "  -0.75%  "

$ws.Range("D36").Value = "0.611"
$ws.Range("E36").Value = "  +1.56%  "

$ws.Range("E37").Value = "  -1.13%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.0165"
$ws.Range("E38").Value = "  -2.78%  "

$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "1.09"
$ws.Range("E39").Value = "  -6.81%  "

$ws.Range("D40").Value = "0.815"
$ws.Range("E40").Value = "  -2.19%  "

$ws.Range("D41").Value = "5.57"
$ws.Range("E41").Value = "  +2.86%  "

$ws.Range("E42").Value = "  -2.15%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.763"
$ws.Range("E43").Value = "  -2.48%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "62.41"
$ws.Range("E44").Value = "  -0.71%  "

$ws.Range("D45").Value = "1.712.22"
$ws.Range("E45").Value = "  -1.41%  "

$ws.Range("D46").Value = "89.17"
$ws.Range("E46").Value = "  -1.45%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0104"
$ws.Range("E47").Value = "  +1.37%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.51"
$ws.Range("E48").Value = "  -2.14%  "

$ws.Range("D49").Value = "0.100"
$ws.Range("E49").Value = "  -1.55%  "

$ws.Range("E50").Value = "  -1.89%  "

$ws.Range("E51").Value = "  -0.23%  "
